$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-26: columns D, J, K, L, M, N, O, P change as new weekly
# data arrived and the consolidated sheet was regenerated/reordered.

$ws.Cells.Item(2, 4).Value = 44162
$ws.Cells.Item(2, 10).Value = 35
$ws.Cells.Item(2, 11).Value = 17000
$ws.Cells.Item(2, 12).Value = 17000
$ws.Cells.Item(2, 13).Value = 17000
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(2, 16).Value = 680

$ws.Cells.Item(3, 4).Value = 44181
$ws.Cells.Item(3, 10).Value = 38
$ws.Cells.Item(3, 11).Value = 26000
$ws.Cells.Item(3, 12).Value = 26000
$ws.Cells.Item(3, 13).Value = 26000
$ws.Cells.Item(3, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 16).Value = 1040

$ws.Cells.Item(4, 4).Value = 44406
$ws.Cells.Item(4, 10).Value = 35
$ws.Cells.Item(4, 11).Value = 32000
$ws.Cells.Item(4, 12).Value = 32000
$ws.Cells.Item(4, 13).Value = 32000
$ws.Cells.Item(4, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(4, 16).Value = 1280

$ws.Cells.Item(5, 4).Value = 44253
$ws.Cells.Item(5, 10).Value = 38
$ws.Cells.Item(5, 11).Value = 18000
$ws.Cells.Item(5, 12).Value = 18000
$ws.Cells.Item(5, 13).Value = 18000
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(5, 16).Value = 720

$ws.Cells.Item(6, 4).Value = 44160
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(6, 11).Value = 21000
$ws.Cells.Item(6, 12).Value = 21000
$ws.Cells.Item(6, 13).Value = 21000
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(6, 16).Value = 840

$ws.Cells.Item(7, 4).Value = 44365
$ws.Cells.Item(7, 10).Value = 70
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 23000
$ws.Cells.Item(7, 13).Value = 22500
$ws.Cells.Item(7, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 16).Value = 900

$ws.Cells.Item(8, 4).Value = 44376
$ws.Cells.Item(8, 10).Value = 38
$ws.Cells.Item(8, 11).Value = 27000
$ws.Cells.Item(8, 12).Value = 27000
$ws.Cells.Item(8, 13).Value = 27000
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(8, 16).Value = 1080

$ws.Cells.Item(9, 4).Value = 44399
$ws.Cells.Item(9, 10).Value = 38
$ws.Cells.Item(9, 11).Value = 33000
$ws.Cells.Item(9, 12).Value = 33000
$ws.Cells.Item(9, 13).Value = 33000
$ws.Cells.Item(9, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(9, 16).Value = 1320

$ws.Cells.Item(10, 4).Value = 44370
$ws.Cells.Item(10, 10).Value = 45
$ws.Cells.Item(10, 11).Value = 32000
$ws.Cells.Item(10, 12).Value = 32000
$ws.Cells.Item(10, 13).Value = 32000
$ws.Cells.Item(10, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(10, 16).Value = 1280

$ws.Cells.Item(11, 4).Value = 44412
$ws.Cells.Item(11, 10).Value = 35
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 24000
$ws.Cells.Item(11, 13).Value = 24000
$ws.Cells.Item(11, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 16).Value = 960

$ws.Cells.Item(12, 4).Value = 44343
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 28000
$ws.Cells.Item(12, 12).Value = 28000
$ws.Cells.Item(12, 13).Value = 28000
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 16).Value = 1120

$ws.Cells.Item(13, 4).Value = 44410
$ws.Cells.Item(13, 10).Value = 35
$ws.Cells.Item(13, 11).Value = 34000
$ws.Cells.Item(13, 12).Value = 34000
$ws.Cells.Item(13, 13).Value = 34000
$ws.Cells.Item(13, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(13, 16).Value = 1360

$ws.Cells.Item(14, 4).Value = 44473
$ws.Cells.Item(14, 10).Value = 85
$ws.Cells.Item(14, 11).Value = 35000
$ws.Cells.Item(14, 12).Value = 36000
$ws.Cells.Item(14, 13).Value = 35471
$ws.Cells.Item(14, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(14, 16).Value = 1419

$ws.Cells.Item(15, 4).Value = 44411
$ws.Cells.Item(15, 10).Value = 35
$ws.Cells.Item(15, 11).Value = 34000
$ws.Cells.Item(15, 12).Value = 34000
$ws.Cells.Item(15, 13).Value = 34000
$ws.Cells.Item(15, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(15, 16).Value = 1360

$ws.Cells.Item(16, 4).Value = 44469
$ws.Cells.Item(16, 10).Value = 73
$ws.Cells.Item(16, 11).Value = 28000
$ws.Cells.Item(16, 12).Value = 29000
$ws.Cells.Item(16, 13).Value = 28521
$ws.Cells.Item(16, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 16).Value = 1141

$ws.Cells.Item(17, 4).Value = 44161
$ws.Cells.Item(17, 10).Value = 35
$ws.Cells.Item(17, 11).Value = 21000
$ws.Cells.Item(17, 12).Value = 21000
$ws.Cells.Item(17, 13).Value = 21000
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(17, 16).Value = 840

$ws.Cells.Item(18, 4).Value = 44165
$ws.Cells.Item(18, 10).Value = 45
$ws.Cells.Item(18, 11).Value = 22000
$ws.Cells.Item(18, 12).Value = 22000
$ws.Cells.Item(18, 13).Value = 22000
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(18, 16).Value = 880

$ws.Cells.Item(19, 4).Value = 44475
$ws.Cells.Item(19, 10).Value = 73
$ws.Cells.Item(19, 11).Value = 25000
$ws.Cells.Item(19, 12).Value = 26000
$ws.Cells.Item(19, 13).Value = 25479
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 16).Value = 1019

$ws.Cells.Item(20, 4).Value = 44252
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 18000
$ws.Cells.Item(20, 12).Value = 19000
$ws.Cells.Item(20, 13).Value = 18625
$ws.Cells.Item(20, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(20, 16).Value = 745

$ws.Cells.Item(21, 4).Value = 44372
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 33000
$ws.Cells.Item(21, 12).Value = 34000
$ws.Cells.Item(21, 13).Value = 33500
$ws.Cells.Item(21, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(21, 16).Value = 1340

$ws.Cells.Item(22, 4).Value = 44448
$ws.Cells.Item(22, 10).Value = 45
$ws.Cells.Item(22, 11).Value = 32000
$ws.Cells.Item(22, 12).Value = 32000
$ws.Cells.Item(22, 13).Value = 32000
$ws.Cells.Item(22, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(22, 16).Value = 1280

$ws.Cells.Item(23, 4).Value = 44452
$ws.Cells.Item(23, 10).Value = 70
$ws.Cells.Item(23, 11).Value = 31000
$ws.Cells.Item(23, 12).Value = 32000
$ws.Cells.Item(23, 13).Value = 31500
$ws.Cells.Item(23, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(23, 16).Value = 1260

$ws.Cells.Item(24, 4).Value = 44453
$ws.Cells.Item(24, 10).Value = 73
$ws.Cells.Item(24, 11).Value = 21000
$ws.Cells.Item(24, 12).Value = 22000
$ws.Cells.Item(24, 13).Value = 21521
$ws.Cells.Item(24, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(24, 16).Value = 861

$ws.Cells.Item(25, 4).Value = 44159
$ws.Cells.Item(25, 10).Value = 35
$ws.Cells.Item(25, 11).Value = 22000
$ws.Cells.Item(25, 12).Value = 22000
$ws.Cells.Item(25, 13).Value = 22000
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(25, 16).Value = 880

$ws.Cells.Item(26, 4).Value = 44468
$ws.Cells.Item(26, 10).Value = 65
$ws.Cells.Item(26, 11).Value = 24000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 24538
$ws.Cells.Item(26, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 16).Value = 982

# Append new row 27 (new weekly record)
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(27, 3).Value = 'Coquimbo'
$ws.Cells.Item(27, 4).Value = 44250
$ws.Cells.Item(27, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = 'Arveja Verde'
$ws.Cells.Item(27, 8).Value = 'Perfection'
$ws.Cells.Item(27, 9).Value = 'Primera'
$ws.Cells.Item(27, 10).Value = 38
$ws.Cells.Item(27, 11).Value = 18000
$ws.Cells.Item(27, 12).Value = 18000
$ws.Cells.Item(27, 13).Value = 18000
$ws.Cells.Item(27, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(27, 16).Value = 720
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = 'Hortaliza'
